{"js": "// Replace the date line and each \"a\u00f7b=c, d\" cell in the worksheet table\n// with the new values from the commit. Every source string in the document\n// is unique, so a scoped, case-sensitive search + Replace is unambiguous.\nconst replacements = [\n  [\"2025-10-17 Friday\", \"2025-10-18 Saturday\"],\n  [\"302\u00f76=50, 2\", \"901\u00f77=128, 5\"],\n  [\"463\u00f72=231, 1\", \"717\u00f72=358, 1\"],\n  [\"977\u00f73=325, 2\", \"892\u00f77=127, 3\"],\n  [\"362\u00f75=72, 2\", \"158\u00f73=52, 2\"],\n  [\"389\u00f77=55, 4\", \"772\u00f72=386, 0\"],\n  [\"999\u00f78=124, 7\", \"457\u00f79=50, 7\"],\n  [\"203\u00f72=101, 1\", \"419\u00f73=139, 2\"],\n  [\"114\u00f75=22, 4\", \"602\u00f77=86, 0\"],\n  [\"801\u00f78=100, 1\", \"512\u00f73=170, 2\"],\n  [\"332\u00f76=55, 2\", \"906\u00f77=129, 3\"],\n  [\"814\u00f79=90, 4\", \"123\u00f72=61, 1\"],\n  [\"795\u00f76=132, 3\", \"978\u00f74=244, 2\"],\n  [\"522\u00f79=58, 0\", \"615\u00f79=68, 3\"],\n  [\"227\u00f76=37, 5\", \"214\u00f73=71, 1\"],\n  [\"316\u00f75=63, 1\", \"173\u00f72=86, 1\"],\n  [\"842\u00f78=105, 2\", \"973\u00f75=194, 3\"],\n  [\"144\u00f77=20, 4\", \"313\u00f74=78, 1\"],\n  [\"592\u00f73=197, 1\", \"284\u00f77=40, 4\"],\n  [\"369\u00f73=123, 0\", \"606\u00f72=303, 0\"],\n  [\"226\u00f74=56, 2\", \"847\u00f78=105, 7\"],\n  [\"492\u00f78=61, 4\", \"908\u00f77=129, 5\"],\n  [\"120\u00f74=30, 0\", \"501\u00f76=83, 3\"],\n  [\"819\u00f77=117, 0\", \"278\u00f76=46, 2\"],\n  [\"775\u00f76=129, 1\", \"198\u00f76=33, 0\"],\n  [\"137\u00f74=34, 1\", \"615\u00f78=76, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const hits = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const hit of hits.items) {\n    hit.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=c, d\" cell in the worksheet table\n# with the new values from the commit. Every source string in the document\n# is unique, so a case-sensitive whole-document Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-10-17 Friday\", \"2025-10-18 Saturday\"),\n  @(\"302\u00f76=50, 2\", \"901\u00f77=128, 5\"),\n  @(\"463\u00f72=231, 1\", \"717\u00f72=358, 1\"),\n  @(\"977\u00f73=325, 2\", \"892\u00f77=127, 3\"),\n  @(\"362\u00f75=72, 2\", \"158\u00f73=52, 2\"),\n  @(\"389\u00f77=55, 4\", \"772\u00f72=386, 0\"),\n  @(\"999\u00f78=124, 7\", \"457\u00f79=50, 7\"),\n  @(\"203\u00f72=101, 1\", \"419\u00f73=139, 2\"),\n  @(\"114\u00f75=22, 4\", \"602\u00f77=86, 0\"),\n  @(\"801\u00f78=100, 1\", \"512\u00f73=170, 2\"),\n  @(\"332\u00f76=55, 2\", \"906\u00f77=129, 3\"),\n  @(\"814\u00f79=90, 4\", \"123\u00f72=61, 1\"),\n  @(\"795\u00f76=132, 3\", \"978\u00f74=244, 2\"),\n  @(\"522\u00f79=58, 0\", \"615\u00f79=68, 3\"),\n  @(\"227\u00f76=37, 5\", \"214\u00f73=71, 1\"),\n  @(\"316\u00f75=63, 1\", \"173\u00f72=86, 1\"),\n  @(\"842\u00f78=105, 2\", \"973\u00f75=194, 3\"),\n  @(\"144\u00f77=20, 4\", \"313\u00f74=78, 1\"),\n  @(\"592\u00f73=197, 1\", \"284\u00f77=40, 4\"),\n  @(\"369\u00f73=123, 0\", \"606\u00f72=303, 0\"),\n  @(\"226\u00f74=56, 2\", \"847\u00f78=105, 7\"),\n  @(\"492\u00f78=61, 4\", \"908\u00f77=129, 5\"),\n  @(\"120\u00f74=30, 0\", \"501\u00f76=83, 3\"),\n  @(\"819\u00f77=117, 0\", \"278\u00f76=46, 2\"),\n  @(\"775\u00f76=129, 1\", \"198\u00f76=33, 0\"),\n  @(\"137\u00f74=34, 1\", \"615\u00f78=76, 7\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n\n$d.Saved = $false\n"}
